$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 31251230
$ws.Range("I137").Value = 43479170
$ws.Range("J137").Value = 2056.111
$ws.Range("K137").Value = 130437510
$ws.Range("L137").Value = 6168.333
$ws.Range("M137").Value = -130434960
$ws.Range("N137").Value = -11268.333

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 60000
$ws.Range("J24").Value = 60000
$ws.Range("L24").Value = 60000
$ws.Range("N24").Value = -60748
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
# Row 32
$ws.Range("H32").Value = 20185.193
$ws.Range("I32").Value = 3559.3103
$ws.Range("K32").Value = 3559.3103
$ws.Range("M32").Value = -3272.3103
# Row 63
$ws.Range("H63").Value = 5806.25
$ws.Range("I63").Value = 5564.2856
$ws.Range("K63").Value = 5564.2856
$ws.Range("M63").Value = -4878.2856
# Row 66
$ws.Range("H66").Value = 5806.25
$ws.Range("I66").Value = 5564.2856
$ws.Range("K66").Value = 27821.428
$ws.Range("M66").Value = -24389.428
# Row 74
$ws.Range("H74").Value = 5767.9644
$ws.Range("I74").Value = 1153.5555
$ws.Range("J74").Value = 14073.9
$ws.Range("K74").Value = 1153.5555
$ws.Range("L74").Value = 14073.9
$ws.Range("M74").Value = -279.5554999999999
$ws.Range("N74").Value = -15821.9
# Row 77
$ws.Range("H77").Value = 5767.9644
$ws.Range("I77").Value = 1153.5555
$ws.Range("J77").Value = 14073.9
$ws.Range("K77").Value = 5767.7775
$ws.Range("L77").Value = 70369.5
$ws.Range("M77").Value = -1399.7775
$ws.Range("N77").Value = -79105.5
# Row 94
$ws.Range("H94").Value = 34658
$ws.Range("J94").Value = 34658
$ws.Range("L94").Value = 34658
$ws.Range("N94").Value = -36460
# Row 100
$ws.Range("H100").Value = 60000
$ws.Range("J100").Value = 60000
$ws.Range("L100").Value = 60000
$ws.Range("N100").Value = -62164
# Row 110
$ws.Range("H110").Value = 1114.4286
$ws.Range("I110").Value = 1050.1666
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1050.1666
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 994.8334
$ws.Range("N110").Value = -5590
# Row 132
$ws.Range("H132").Value = 2650.125
$ws.Range("I132").Value = 2251.2334
$ws.Range("J132").Value = 3846.8
$ws.Range("K132").Value = 6753.7002
$ws.Range("L132").Value = 11540.4
$ws.Range("M132").Value = -4223.7002
$ws.Range("N132").Value = -16600.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1101.125
$ws.Range("I94").Value = 1184.8334
$ws.Range("J94").Value = 850
$ws.Range("K94").Value = 1184.8334
$ws.Range("L94").Value = 850
$ws.Range("M94").Value = -733.8334
$ws.Range("N94").Value = -1752
# Row 105
$ws.Range("H105").Value = 288785.78
$ws.Range("I105").Value = 2841.4285
$ws.Range("J105").Value = 717702.3
$ws.Range("K105").Value = 2841.4285
$ws.Range("L105").Value = 717702.3
$ws.Range("M105").Value = -1094.4285
$ws.Range("N105").Value = -721196.3
# Row 134
$ws.Range("H134").Value = 3029.8708
$ws.Range("I134").Value = 2116
$ws.Range("J134").Value = 5263.778
$ws.Range("K134").Value = 6348
$ws.Range("L134").Value = 15791.334
$ws.Range("M134").Value = -3813
$ws.Range("N134").Value = -20861.334

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6315
$ws.Range("I31").Value = 2537.1667
$ws.Range("J31").Value = 8701
$ws.Range("K31").Value = 2537.1667
$ws.Range("L31").Value = 8701
$ws.Range("M31").Value = -2242.1667
$ws.Range("N31").Value = -9291
# Row 34
$ws.Range("H34").Value = 6315
$ws.Range("I34").Value = 2537.1667
$ws.Range("J34").Value = 8701
$ws.Range("K34").Value = 2537.1667
$ws.Range("L34").Value = 8701
$ws.Range("M34").Value = -2335.1667
$ws.Range("N34").Value = -9105
# Row 70
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28630
# Row 73
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -30184
# Row 132
$ws.Range("H132").Value = 1993.1774
$ws.Range("I132").Value = 1733.0385
$ws.Range("K132").Value = 5199.1155
$ws.Range("M132").Value = -2669.1155

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 164371.28
$ws.Range("I4").Value = 153938.31
$ws.Range("J4").Value = 300000
$ws.Range("K4").Value = 461814.93
$ws.Range("L4").Value = 900000
$ws.Range("M4").Value = -461702.93
$ws.Range("N4").Value = -900224
# Row 34
$ws.Range("H34").Value = 1725.125
$ws.Range("I34").Value = 4502
$ws.Range("J34").Value = 1328.4286
$ws.Range("K34").Value = 13506
$ws.Range("L34").Value = 3985.2858
$ws.Range("M34").Value = -13422
$ws.Range("N34").Value = -4153.2858
# Row 39
$ws.Range("H39").Value = 9500
$ws.Range("J39").Value = 9500
$ws.Range("L39").Value = 28500
$ws.Range("N39").Value = -29088
# Row 55
$ws.Range("H55").Value = 4000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -12354
# Row 68
$ws.Range("H68").Value = 20726.166
$ws.Range("I68").Value = 24760.4
$ws.Range("J68").Value = 555
$ws.Range("K68").Value = 74281.20000000001
$ws.Range("L68").Value = 1665
$ws.Range("M68").Value = -73470.20000000001
$ws.Range("N68").Value = -3287
# Row 71
$ws.Range("H71").Value = 20726.166
$ws.Range("I71").Value = 24760.4
$ws.Range("J71").Value = 555
$ws.Range("K71").Value = 222843.6
$ws.Range("L71").Value = 4995
$ws.Range("M71").Value = -218787.6
$ws.Range("N71").Value = -13107
# Row 132
$ws.Range("H132").Value = 864.5294
$ws.Range("I132").Value = 782.1667
$ws.Range("J132").Value = 909.4545000000001
$ws.Range("K132").Value = 7039.5003
$ws.Range("L132").Value = 8185.0905
$ws.Range("M132").Value = -4509.5003
$ws.Range("N132").Value = -13245.0905

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2772.375
$ws.Range("I132").Value = 2767
$ws.Range("J132").Value = 2785.4285
$ws.Range("K132").Value = 8301
$ws.Range("L132").Value = 8356.2855
$ws.Range("M132").Value = -5771
$ws.Range("N132").Value = -13416.2855

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 46
$ws.Range("H46").Value = 1105.0769
$ws.Range("I46").Value = 911.3333
$ws.Range("J46").Value = 1271.1428
$ws.Range("K46").Value = 911.3333
$ws.Range("L46").Value = 1271.1428
$ws.Range("M46").Value = -723.3333
$ws.Range("N46").Value = -1647.1428
# Row 132
$ws.Range("H132").Value = 2601.111
$ws.Range("I132").Value = 1734.8518
$ws.Range("J132").Value = 5199.8887
$ws.Range("K132").Value = 5204.555399999999
$ws.Range("L132").Value = 15599.6661
$ws.Range("M132").Value = -2674.555399999999
$ws.Range("N132").Value = -20659.6661

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 44473444
$ws.Range("J2").Value = 52752.25
$ws.Range("L2").Value = 52752.25
$ws.Range("N2").Value = -52976.25
# Row 113
$ws.Range("H113").Value = 283.47058
$ws.Range("I113").Value = 283.47058
$ws.Range("K113").Value = 850.41174
$ws.Range("M113").Value = 1319.58826
# Row 132
$ws.Range("H132").Value = 1653.3469
$ws.Range("I132").Value = 1418.4474
$ws.Range("K132").Value = 4255.3422
$ws.Range("M132").Value = -1725.3422
